$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "Fri Mar 08 00:15:36 EST 2024",
    "Fri Mar 08 00:16:07 EST 2024",
    "Fri Mar 08 00:16:33 EST 2024",
    "Fri Mar 08 00:16:54 EST 2024",
    "Fri Mar 08 00:17:31 EST 2024",
    "Fri Mar 08 00:17:56 EST 2024",
    "Fri Mar 08 00:18:20 EST 2024",
    "Fri Mar 08 00:18:46 EST 2024",
    "Fri Mar 08 00:19:24 EST 2024",
    "Fri Mar 08 00:19:52 EST 2024",
    "Fri Mar 08 00:20:08 EST 2024",
    "Fri Mar 08 00:20:53 EST 2024",
    "Fri Mar 08 00:21:33 EST 2024",
    "Fri Mar 08 00:22:10 EST 2024",
    "Fri Mar 08 00:22:48 EST 2024",
    "Fri Mar 08 00:23:19 EST 2024",
    "Fri Mar 08 00:23:53 EST 2024",
    "Fri Mar 08 00:24:09 EST 2024",
    "Fri Mar 08 00:25:12 EST 2024",
    "Fri Mar 08 00:25:37 EST 2024",
    "Fri Mar 08 00:25:57 EST 2024",
    "Fri Mar 08 00:26:18 EST 2024",
    "Fri Mar 08 00:26:49 EST 2024",
    "Fri Mar 08 00:27:11 EST 2024",
    "Fri Mar 08 00:27:48 EST 2024",
    "Fri Mar 08 00:28:12 EST 2024",
    "Fri Mar 08 00:28:49 EST 2024",
    "Fri Mar 08 00:29:25 EST 2024",
    "Fri Mar 08 00:29:54 EST 2024",
    "Fri Mar 08 00:30:32 EST 2024",
    "Fri Mar 08 00:31:13 EST 2024",
    "Fri Mar 08 00:31:39 EST 2024",
    "Fri Mar 08 00:32:13 EST 2024",
    "Fri Mar 08 00:32:41 EST 2024",
    "Fri Mar 08 00:33:15 EST 2024",
    "Fri Mar 08 00:33:40 EST 2024",
    "Fri Mar 08 00:34:22 EST 2024",
    "Fri Mar 08 00:34:50 EST 2024",
    "Fri Mar 08 00:35:14 EST 2024",
    "Fri Mar 08 00:35:44 EST 2024",
    "Fri Mar 08 00:36:10 EST 2024",
    "Fri Mar 08 00:36:45 EST 2024",
    "Fri Mar 08 00:37:11 EST 2024",
    "Fri Mar 08 00:37:53 EST 2024",
    "Fri Mar 08 00:38:37 EST 2024",
    "Fri Mar 08 00:38:54 EST 2024"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
